# Apply updated evaluation metrics to the workbook

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Summary ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.599250936329588
$summary.Range("C2").Value = 0.5550935550935551
$summary.Range("D2").Value = 1
$summary.Range("E2").Value = 0.713903743315508
$summary.Range("F2").Value = 0.8618463524854745
$summary.Range("G2").Value = 0.9700950251537171
$summary.Range("H2").Value = 0.679494732707711
$summary.Range("I2").Value = 534
$summary.Range("J2").Value = 428
$summary.Range("K2").Value = 106
$summary.Range("L2").Value = 0

# ---- Sheet 2: Classification Report ----
$report = $wb.Worksheets.Item("Classification Report")

$report.Range("B2").Value = 1
$report.Range("C2").Value = 0.198501872659176
$report.Range("D2").Value = 0.33125

$report.Range("B3").Value = 0.5550935550935551
$report.Range("C3").Value = 1
$report.Range("D3").Value = 0.713903743315508

$report.Range("B4").Value = 0.599250936329588
$report.Range("C4").Value = 0.599250936329588
$report.Range("D4").Value = 0.599250936329588
$report.Range("E4").Value = 0.599250936329588

$report.Range("B5").Value = 0.7775467775467775
$report.Range("C5").Value = 0.599250936329588
$report.Range("D5").Value = 0.522576871657754

$report.Range("B6").Value = 0.7775467775467776
$report.Range("C6").Value = 0.599250936329588
$report.Range("D6").Value = 0.522576871657754

# ---- Sheet 3: Confusion Matrix ----
$confusion = $wb.Worksheets.Item("Confusion Matrix")

$confusion.Range("B2").Value = 106
$confusion.Range("C2").Value = 428

$confusion.Range("B3").Value = 0
$confusion.Range("C3").Value = 534
